# normalization-for-heterogeneous-effects.xlsx — update ban relative-risk
# inputs, drop the now-unused "1-p_severe" helper row, and add the
# overall injury-probability check in B17.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- updated input assumptions (B4:B6) ---------------------------------
$ws.Range("B4").Value = 0.8   # rr_severe_ban
$ws.Range("B5").Value = 0.7   # rr_mildmod_ban
$ws.Range("B6").Value = 0.9   # rr_fatality_ban

# --- drop the old "p_mildmod_ban = 1-p_severe" scratch row -------------
$ws.Range("A8:C8").ClearContents()

# --- new summary check: total normalized probability / p_injury --------
$ws.Range("B17").Formula = "=B13/p_injury"

# --- restore the view state (scroll position + active cell) -----------
$ws.Range("A2").Select()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A11").Select()
